$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated summary statistics after switching the removals approach to
# "removal_ha (sem -)" — negative hectare values are no longer kept as
# negatives, so the Sum/Mean/Median/SD/SSE/Min/Max/Skw/Krt figures change
# (Min collapses to 0 since no negative values remain).

$ws.Range("B2").Value = 5353.332271964373
$ws.Range("C2").Value = 0.2011321112099629
$ws.Range("D2").Value = 0.1232157853975029
$ws.Range("E2").Value = 0.354523859628911
$ws.Range("F2").Value = 0.002153738030818717
$ws.Range("G2").Value = -0
$ws.Range("I2").Value = 7.871819021107108
$ws.Range("J2").Value = 85.74701605258493

$ws.Range("B3").Value = 27687.48800413157
$ws.Range("C3").Value = 0.439483936573517
$ws.Range("D3").Value = 0.2884945493498523
$ws.Range("E3").Value = 0.8131875277778049
$ws.Range("F3").Value = 0.003190958764847861
$ws.Range("G3").Value = -0
$ws.Range("I3").Value = 13.88576153367213
$ws.Range("J3").Value = 282.1142438711866

$ws.Range("B4").Value = 8066.396757711413
$ws.Range("C4").Value = 0.1599333166331869
$ws.Range("D4").Value = 0.08873628585163773
$ws.Range("E4").Value = 0.3704295059340151
$ws.Range("F4").Value = 0.001635486222070736
$ws.Range("G4").Value = -0
$ws.Range("I4").Value = 14.13739385690647
$ws.Range("J4").Value = 271.8226083234179
